$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price observation: insert a new row at 220 (existing rows
# 220-262 shift down to 221-263) and populate it with the latest reading.
$ws.Rows.Item(220).Insert()

$ws.Cells.Item(220, 1).Value = 4
$ws.Cells.Item(220, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(220, 3).Value = "Los Lagos"
$ws.Cells.Item(220, 4).Value2 = 44694
$ws.Cells.Item(220, 5).Value = 10
$ws.Cells.Item(220, 6).Value = 100112003
$ws.Cells.Item(220, 7).Value = "Ajo"
$ws.Cells.Item(220, 8).Value = "Chino"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 220
$ws.Cells.Item(220, 11).Value = 21000
$ws.Cells.Item(220, 12).Value = 21000
$ws.Cells.Item(220, 13).Value = 21000
$ws.Cells.Item(220, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(220, 15).Value = "China"
$ws.Cells.Item(220, 16).Value = 2100
$ws.Cells.Item(220, 17).Value = 10
$ws.Cells.Item(220, 18).Value = "Hortaliza"
